$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.133.28"
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("D3").Value = "3.772.32"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "406.38"
$ws.Range("E5").Value = "  -4.23%  "
$ws.Range("D6").Value = "132.20"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "3.764.77"
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").Value = "  -5.95%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "0.721"
$ws.Range("E10").Value = "  -6.13%  "
$ws.Range("D11").Value = "0.166"
$ws.Range("E11").Value = "  -9.27%  "
$ws.Range("D12").Value = "0.0000353"
$ws.Range("E12").Value = "  -9.17%  "
$ws.Range("D13").Value = "40.49"
$ws.Range("E13").Value = "  -5.87%  "
$ws.Range("D14").Value = "4.365.00"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("D15").Value = "9.70"
$ws.Range("E15").Value = "  -5.22%  "
$ws.Range("D16").Value = "14.53"
$ws.Range("E16").Value = "  +11.66%  "
$ws.Range("D17").Value = "0.138"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "3.800.16"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").Value = "19.41"
$ws.Range("E19").Value = "  -7.12%  "
$ws.Range("D20").Value = "66.643.71"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").Value = "1.06"
$ws.Range("E21").Value = "  -6.75%  "
$ws.Range("D22").Value = "412.58"
$ws.Range("E22").Value = "  -8.29%  "
$ws.Range("E23").Value = "  -11.22%  "
$ws.Range("D24").Value = "84.89"
$ws.Range("E24").Value = "  -5.53%  "
$ws.Range("D25").Value = "3.02"
$ws.Range("E25").Value = "  -5.32%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "36.20"
$ws.Range("E26").Value = "  -5.25%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "5.68"
$ws.Range("E27").Value = "  +14.33%  "
$ws.Range("D28").Value = "3.09"
$ws.Range("E28").Value = "  -6.84%  "
$ws.Range("D29").Value = "9.29"
$ws.Range("E29").Value = "  -9.30%  "
$ws.Range("D30").Value = "12.33"
$ws.Range("E30").Value = "  -3.26%  "
$ws.Range("D31").Value = "2.73"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").Value = "0.118"
$ws.Range("E32").Value = "  -4.86%  "
$ws.Range("D33").Value = "7.17"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("E34").Value = "  -6.37%  "
$ws.Range("D35").Value = "38.93"
$ws.Range("E35").Value = "  -8.49%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "55.12"
$ws.Range("E37").Value = "  -2.34%  "
$ws.Range("D38").Value = "0.0₃0741"
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("D39").Value = "0.0457"
$ws.Range("E39").Value = "  -7.31%  "
$ws.Range("D40").Value = "2.85"
$ws.Range("E40").Value = "  -7.92%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.135"
$ws.Range("E42").Value = "  -8.33%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "26.82"
$ws.Range("E43").Value = "  -6.56%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "144.68"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.08"
$ws.Range("E45").Value = "  +16.47%  "
$ws.Range("D46").Value = "3.24"
$ws.Range("E46").Value = "  -6.80%  "
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("E48").Value = "  -5.91%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "2.57"
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "4.22"
$ws.Range("E50").Value = "  -5.24%  "
$ws.Range("D51").Value = "0.290"
$ws.Range("E51").Value = "  -7.03%  "